$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Samples" tab query (row 3, column B: TabQuery for SamplesTab) drops the
# two extra columns (Tumor, Analyte Type) that used to be selected, leaving
# just Sample ID / Participant ID / Study Name / Accession.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
  s.phs_accession = 'phs001524' AND gi.library_selection = 'Hybrid Selection'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newSamplesQuery

# Move the active selection to C3 (matches the saved view state).
$ws.Range("C3").Select() | Out-Null
